$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("D9").Value = "1399-03-13 (9)"
$ws.Range("E9").Value = "1400-03-11 (9)"
$ws.Range("F9").Value = "1401-03-11 (10)"
$ws.Range("G9").Value = "1402-02-27 (7)"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-27"
$ws.Range("D11").Value = 1155346
$ws.Range("E11").Value = 1736841
$ws.Range("F11").Value = 2440996
$ws.Range("G11").Value = 4297311
$ws.Range("H11").Value = 9819805
$ws.Range("D12").Value = -794103
$ws.Range("E12").Value = -1061828
$ws.Range("F12").Value = -1200366
$ws.Range("G12").Value = -2127640
$ws.Range("H12").Value = -5906404
$ws.Range("D13").Value = 361243
$ws.Range("E13").Value = 675013
$ws.Range("F13").Value = 1240630
$ws.Range("G13").Value = 2169671
$ws.Range("H13").Value = 3913401
$ws.Range("D14").Value = -51624
$ws.Range("E14").Value = -71156
$ws.Range("F14").Value = -116636
$ws.Range("G14").Value = -177907
$ws.Range("H14").Value = -213869
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("D16").Value = 9872
$ws.Range("E16").Value = 24621
$ws.Range("F16").Value = 15630
$ws.Range("G16").Value = 40427
$ws.Range("H16").Value = 28534
$ws.Range("D17").Value = 319491
$ws.Range("E17").Value = 628478
$ws.Range("F17").Value = 1139624
$ws.Range("G17").Value = 2032191
$ws.Range("H17").Value = 3728066
$ws.Range("D18").Value = -99694
$ws.Range("E18").Value = -156780
$ws.Range("F18").Value = -155571
$ws.Range("G18").Value = -232392
$ws.Range("H18").Value = -337504
$ws.Range("D19").Value = -5559
$ws.Range("E19").Value = -121652
$ws.Range("F19").Value = -69982
$ws.Range("G19").Value = -33612
$ws.Range("H19").Value = -2354
$ws.Range("D20").Value = 214238
$ws.Range("E20").Value = 350046
$ws.Range("F20").Value = 914071
$ws.Range("G20").Value = 1766187
$ws.Range("H20").Value = 3388208
$ws.Range("D21").Value = -52867
$ws.Range("E21").Value = -66781
$ws.Range("F21").Value = -110014
$ws.Range("G21").Value = -248862
$ws.Range("H21").Value = -501836
$ws.Range("D22").Value = 161371
$ws.Range("E22").Value = 283265
$ws.Range("F22").Value = 804057
$ws.Range("G22").Value = 1517325
$ws.Range("H22").Value = 2886372
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("D24").Value = 161371
$ws.Range("E24").Value = 283265
$ws.Range("F24").Value = 804057
$ws.Range("G24").Value = 1517325
$ws.Range("H24").Value = 2886372
$ws.Range("D25").Value = 633
$ws.Range("E25").Value = 1111
$ws.Range("F25").Value = 3153
$ws.Range("G25").Value = 3335
$ws.Range("H25").Value = 3207
$ws.Range("D26").Value = 255000
$ws.Range("E26").Value = 255000
$ws.Range("F26").Value = 255000
$ws.Range("G26").Value = 455000
$ws.Range("H26").Value = 900000
$ws.Range("D27").Value = 179
$ws.Range("E27").Value = 315
$ws.Range("F27").Value = 893
$ws.Range("G27").Value = 1686
$ws.Range("H27").Value = 3207
